# Apply the three title-text edits (add the "STAMP " prefix) and the
# Handout Master date placeholder update (2/23/21 -> 2/24/21).

$p = $ppt.ActivePresentation

# --- Slide titles: prepend "STAMP " to the existing title text ---------

$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "STAMP Return Path Control Code Sub-TLV - Usage"

$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "STAMP Return Address Sub-TLV - Usage"

$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "STAMP Return Path Segment List Sub-TLVs - Usage"

# --- Handout Master "date" placeholder (auto date field) ----------------

$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$hf.DateAndTime.Text = "2/24/21"
